$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "aviv12"
$ws.Range("B10").Value = "aviv123!"
